$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header info updates
$ws.Range("C2").Value = "Hartmut"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Opening balance date
$ws.Range("D5").Value = "KONTOSTAND AM 03.04.2024"

# Row 6
$ws.Range("B6").Value = "06.04."
$ws.Range("C6").Value = "07.04."
$ws.Range("D6").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 68597295"
$ws.Range("E6").Value = "84,08-"

# Row 7
$ws.Range("B7").Value = "08.04."
$ws.Range("C7").Value = "09.04."
$ws.Range("D7").Value = "KARTENZAHLUNG SHELL TANKSTELLE"
$ws.Range("E7").Value = "49,79-"

# Row 8
$ws.Range("B8").Value = "12.04."
$ws.Range("C8").Value = "13.04."
$ws.Range("D8").Value = "PAYPAL EDQVHM"
$ws.Range("E8").Value = "99,98-"

# Row 9 (previously empty) - copy the amount-cell number format/alignment (s=17) from row 6's E cell
$ws.Range("B9").Value = "15.04."
$ws.Range("C9").Value = "16.04."
$ws.Range("D9").Value = "AMAZON.DE MKTPLC EU ZLWRIT"
$ws.Range("E6").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("E9").Value = "163,94-"

# Row 10 (previously empty)
$ws.Range("B10").Value = "16.04."
$ws.Range("C10").Value = "17.04."
$ws.Range("D10").Value = "PAYPAL MLDIXC"
$ws.Range("E6").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("E10").Value = "93,96-"

# Row 11 (previously empty)
$ws.Range("B11").Value = "19.04."
$ws.Range("C11").Value = "20.04."
$ws.Range("D11").Value = "KARTENZ./19.04 REWE RO"
$ws.Range("E6").Copy()
$ws.Range("E11").PasteSpecial(-4122)
$ws.Range("E11").Value = "143,96-"

# Closing balance
$ws.Range("D12").Value = "KONTOSTAND AM 21.04.2024"
$ws.Range("E12").Value = "635,71-"

# Next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 26.04.2024"
